$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Batch_sql id labels for the existing batch_010/batch_011 cases
# (they were renamed from two-digit to three-digit suffixes: _10 -> _010, _11 -> _011)
$ws.Range("I11").Value = "batch_sql_010"
$ws.Range("I12").Value = "batch_sql_011"

# Add a new test case row for batch_012 (one more case about auto_increment)
$ws.Range("A13").Value = "batch_012"
$ws.Range("B13").Value = "y"
$ws.Range("C13").Value = "批量操作语句12执行"
$ws.Range("D13").Value = "batchsql"
$ws.Range("E13").Value = "SingleTable"
$ws.Range("G13").Value = "batch012"
$ws.Range("I13").Value = "batch_sql_012"
$ws.Range("J13").Value = "select * from `$batch012"
$ws.Range("K13").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/batchsql/expectedresult/batch_012.csv"
$ws.Range("N13").Value = "csv_containsAll"

# Match the "fill" alignment formatting used by the other rows in the
# Query_result1 (K) column
$ws.Range("K12").Copy()
$ws.Range("K13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active cell selection to match the edited workbook
[void]$ws.Range("E22").Select()
